{"js": "const replacements = [\n  [\"2025-10-29 Wednesday\", \"2025-10-30 Thursday\"],\n  [\"111\u00d79=\", \"171\u00d72=\"],\n  [\"458\u00d74=\", \"151\u00d78=\"],\n  [\"941\u00d73=\", \"981\u00d73=\"],\n  [\"572\u00d74=\", \"678\u00d79=\"],\n  [\"242\u00d79=\", \"245\u00d77=\"],\n  [\"814\u00d79=\", \"465\u00d75=\"],\n  [\"353\u00d78=\", \"848\u00d76=\"],\n  [\"169\u00d72=\", \"602\u00d73=\"],\n  [\"668\u00d79=\", \"581\u00d74=\"],\n  [\"795\u00d79=\", \"365\u00d77=\"],\n  [\"433\u00d76=\", \"709\u00d74=\"],\n  [\"841\u00d75=\", \"299\u00d76=\"],\n  [\"834\u00d72=\", \"525\u00d79=\"],\n  [\"118\u00d74=\", \"340\u00d77=\"],\n  [\"543\u00d77=\", \"699\u00d78=\"],\n  [\"539\u00d79=\", \"554\u00d75=\"],\n  [\"802\u00d78=\", \"507\u00d78=\"],\n  [\"145\u00d74=\", \"693\u00d79=\"],\n  [\"256\u00d78=\", \"602\u00d73=\"],\n  [\"653\u00d74=\", \"443\u00d79=\"],\n  [\"581\u00d73=\", \"837\u00d74=\"],\n  [\"390\u00d73=\", \"720\u00d78=\"],\n  [\"970\u00d78=\", \"788\u00d73=\"],\n  [\"753\u00d75=\", \"414\u00d73=\"],\n  [\"267\u00d72=\", \"231\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-29 Wednesday\", \"2025-10-30 Thursday\"),\n    @(\"111\u00d79=\", \"171\u00d72=\"),\n    @(\"458\u00d74=\", \"151\u00d78=\"),\n    @(\"941\u00d73=\", \"981\u00d73=\"),\n    @(\"572\u00d74=\", \"678\u00d79=\"),\n    @(\"242\u00d79=\", \"245\u00d77=\"),\n    @(\"814\u00d79=\", \"465\u00d75=\"),\n    @(\"353\u00d78=\", \"848\u00d76=\"),\n    @(\"169\u00d72=\", \"602\u00d73=\"),\n    @(\"668\u00d79=\", \"581\u00d74=\"),\n    @(\"795\u00d79=\", \"365\u00d77=\"),\n    @(\"433\u00d76=\", \"709\u00d74=\"),\n    @(\"841\u00d75=\", \"299\u00d76=\"),\n    @(\"834\u00d72=\", \"525\u00d79=\"),\n    @(\"118\u00d74=\", \"340\u00d77=\"),\n    @(\"543\u00d77=\", \"699\u00d78=\"),\n    @(\"539\u00d79=\", \"554\u00d75=\"),\n    @(\"802\u00d78=\", \"507\u00d78=\"),\n    @(\"145\u00d74=\", \"693\u00d79=\"),\n    @(\"256\u00d78=\", \"602\u00d73=\"),\n    @(\"653\u00d74=\", \"443\u00d79=\"),\n    @(\"581\u00d73=\", \"837\u00d74=\"),\n    @(\"390\u00d73=\", \"720\u00d78=\"),\n    @(\"970\u00d78=\", \"788\u00d73=\"),\n    @(\"753\u00d75=\", \"414\u00d73=\"),\n    @(\"267\u00d72=\", \"231\u00d74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute(\n        $find.Text,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $true,\n        $find.Replacement.Text,\n        2\n    ) | Out-Null\n}\n"}
